$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row (period 2509) right after the last existing
#     worker/period row (row 17), pushing the blank gap + signature block
#     down by one row. ---
$ws.Rows.Item(18).Insert()

# Copy formatting from the row directly above (row 17) so the new row
# matches the existing table styling (borders/fill/font).
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# New row values: same worker, new period (2509), same amounts as the
# other rows for this worker.
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143127438"
$ws.Range("D18").Value = "LUIS FERNANDO MIRANDA GUTIERREZ"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2509"
$ws.Range("F18").Value = 166480
$ws.Range("G18").Value = 4162000

# --- Update the summary figures on the account-statement header. ---
# Valor Mora (total overdue amount) grew with the new period.
$ws.Range("E11").Value = 499440

# Cant. Periodos (number of overdue periods) grew from 2 to 3.
$ws.Range("F13").Value = 3
